$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the A84 timestamp (tiny floating point precision fix)
$ws.Range("A84").Value = 44397.76866720023

# Append the newly retrieved row of data as row 85
$ws.Range("A85").Value = 44398.7692034377
$ws.Range("B85").Value = 80186
$ws.Range("C85").Value = 67655
$ws.Range("D85").Value = 3758
$ws.Range("E85").Value = 2196
$ws.Range("F85").Value = 1590
$ws.Range("G85").Value = 20969
$ws.Range("H85").Value = 1625
$ws.Range("I85").Value = 888
$ws.Range("J85").Value = 198
